# No-op: the source diff for this fixture consists solely of XML attribute /
# namespace-declaration re-ordering inside ppt/slideMasters/slideMaster1.xml
# (e.g. <p:ph type="body" idx="1"/> -> <p:ph idx="1" type="body"/>), a
# byte-level serialization artifact of the original tooling's canonical
# (alphabetically-sorted) XML writer. Every attribute name/value pair is
# identical before and after; nothing in the slide master's content,
# formatting, placeholders, text styles, color map, or theme references
# actually changes. There is no PowerPoint object-model action that
# corresponds to "re-order XML attributes", so there is nothing to apply
# here -- the presentation is left exactly as loaded.
$p = $ppt.ActivePresentation
